# Update Betfair Back/Lay odds values for 2025-11-10 workbook
# (refreshed odds snapshot across several rows/columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.74
$ws.Range("G2").Value = 1.85
$ws.Range("H2").Value = 4.8
$ws.Range("I2").Value = 5.7
$ws.Range("J2").Value = 3.85
$ws.Range("L2").Value = 1.32
$ws.Range("N2").Value = 3.75
$ws.Range("O2").Value = 1.3
$ws.Range("P2").Value = 1.95
$ws.Range("Q2").Value = 1.84
$ws.Range("R2").Value = 1.38
$ws.Range("S2").Value = 3.15
$ws.Range("T2").Value = 1.79
$ws.Range("V2").Value = 1.21
$ws.Range("W2").Value = 2.16
$ws.Range("AE2").Value = 75
$ws.Range("AM2").Value = 120
$ws.Range("AN2").Value = 13.5
$ws.Range("F4").Value = 1.72
$ws.Range("K4").Value = 5
$ws.Range("P4").Value = 2.12
$ws.Range("Q4").Value = 1.73
$ws.Range("T4").Value = 1.74
$ws.Range("U4").Value = 2.14
$ws.Range("W4").Value = 2.16
$ws.Range("AE4").Value = 65
$ws.Range("X5").Value = 46
$ws.Range("Y5").Value = 29
$ws.Range("Z5").Value = 34
$ws.Range("AB5").Value = 26
$ws.Range("AC5").Value = 970
$ws.Range("AD5").Value = 970
$ws.Range("AE5").Value = 32
$ws.Range("AF5").Value = 28
$ws.Range("AG5").Value = 970
$ws.Range("AH5").Value = 970
$ws.Range("AI5").Value = 34
$ws.Range("AJ5").Value = 40
$ws.Range("AK5").Value = 25
$ws.Range("AL5").Value = 29
$ws.Range("AM5").Value = 50
$ws.Range("AN5").Value = 970
$ws.Range("AO5").Value = 970
$ws.Range("F6").Value = 2.18
$ws.Range("G6").Value = 2.54
$ws.Range("H6").Value = 3.2
$ws.Range("I6").Value = 3.9
$ws.Range("J6").Value = 3.1
$ws.Range("K6").Value = 3.8
$ws.Range("O6").Value = 1.35
$ws.Range("Q6").Value = 2.02
$ws.Range("S6").Value = 3.1
$ws.Range("T6").Value = 1.79
$ws.Range("U6").Value = 1.96
$ws.Range("V6").Value = 1.35
$ws.Range("W6").Value = 1.65
$ws.Range("Y6").Value = 13
$ws.Range("Z6").Value = 26
$ws.Range("AA6").Value = 75
$ws.Range("AD6").Value = 16
$ws.Range("AE6").Value = 48
$ws.Range("AF6").Value = 16
$ws.Range("AI6").Value = 60
$ws.Range("AJ6").Value = 36
$ws.Range("AK6").Value = 30
$ws.Range("AO6").Value = 55
$ws.Range("F7").Value = 2.6
$ws.Range("G7").Value = 2.8
$ws.Range("H7").Value = 2.96
$ws.Range("I7").Value = 3.25
$ws.Range("V7").Value = 1.45
$ws.Range("W7").Value = 1.55
$ws.Range("AB7").Value = 970
$ws.Range("AD7").Value = 970
$ws.Range("AE7").Value = 40
$ws.Range("AF7").Value = 970
$ws.Range("AG7").Value = 970
$ws.Range("AJ7").Value = 44
$ws.Range("AK7").Value = 34
$ws.Range("AN7").Value = 32
$ws.Range("AO7").Value = 46
$ws.Range("J8").Value = 2.82
$ws.Range("N8").Value = 2.3
$ws.Range("P8").Value = 1.43
$ws.Range("U8").Value = 1.66
$ws.Range("V8").Value = 1.63
$ws.Range("AM8").Value = 320
$ws.Range("H9").Value = 1.11
$ws.Range("J9").Value = 1.13
$ws.Range("O9").Value = 1.25
$ws.Range("P9").Value = 1.25
$ws.Range("Q9").Value = 1.85
$ws.Range("V9").Value = 1.3
$ws.Range("W9").Value = 1.18
$ws.Range("Z11").Value = 32
$ws.Range("AJ11").Value = 28
$ws.Range("AL11").Value = 48
$ws.Range("V13").Value = 1.13
$ws.Range("X13").Value = 960
$ws.Range("Y13").Value = 28
$ws.Range("AB13").Value = 970
$ws.Range("AC13").Value = 970
$ws.Range("AD13").Value = 34
$ws.Range("AF13").Value = 970
$ws.Range("AG13").Value = 970
$ws.Range("AH13").Value = 29
$ws.Range("AJ13").Value = 970
$ws.Range("AK13").Value = 960
$ws.Range("AL13").Value = 46
$ws.Range("F14").Value = 2.54
$ws.Range("G14").Value = 2.8
$ws.Range("I14").Value = 3.55
$ws.Range("V14").Value = 1.4
$ws.Range("W14").Value = 1.55
$ws.Range("Y14").Value = 980
